$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.823.37"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.740.42"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.44"
$ws.Range("E5").Value = "  -4.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5147"
$ws.Range("E7").Value = "  +1.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2740"
$ws.Range("E8").Value = "  +2.94%  "
$ws.Range("E9").Value = "  -5.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06087"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").Value = "1.740.31"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07000"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.13"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6309"
$ws.Range("E14").Value = "  +4.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.495"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.33"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "25.837.36"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.44"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006607"
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("D22").Value = "1.958.15"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.073"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.434"
$ws.Range("E24").Value = "  +3.02%  "
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.22"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.492"
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.807"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.92"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.88"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08267"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.596"
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.368"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04384"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9627"
$ws.Range("E36").Value = "  -3.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5950"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.660"
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01546"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.923"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.71"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3796"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7280"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.865"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05476"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.234"
$ws.Range("E47").Value = "  +5.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1099"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.57"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.91"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.458"
$ws.Range("E51").Value = "  -2.48%  "
